$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("user_00")

$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 0
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 0
